$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.552.03'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.445.75'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.69'
$ws.Range("E5").Value = '  -2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.21'
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.443.63'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("E12").Value = '  -2.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.040.57'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.94'
$ws.Range("E14").Value = '  -3.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.132'
$ws.Range("E15").Value = '  -3.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.575.00'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.443.46'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("E19").Value = '  -3.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.81'
$ws.Range("E20").Value = '  -3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.89'
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.69'
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.72'
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.526'
$ws.Range("E26").Value = '  -1.60%  '
$ws.Range("E27").Value = '  -2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.83'
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  -5.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.90'
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("E34").Value = '  -5.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.03'
$ws.Range("E36").Value = '  -4.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.60'
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.876'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.11'
$ws.Range("E40").Value = '  +4.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  -4.97%  '
$ws.Range("E42").Value = '  -4.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.53'
$ws.Range("E43").Value = '  -5.08%  '
$ws.Range("E44").Value = '  -3.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.693.07'
$ws.Range("E45").Value = '  -4.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0692'
$ws.Range("E46").Value = '  -4.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.18'
$ws.Range("E47").Value = '  -5.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.35'
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '320.24'
$ws.Range("E50").Value = '  -4.71%  '
$ws.Range("E51").Value = '  -3.99%  '
